$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1063.3334
$ws.Range("I70").Value = 1193.3334
$ws.Range("J70").Value = 933.3333
$ws.Range("K70").Value = 3580.0002
$ws.Range("L70").Value = 2799.9999
$ws.Range("M70").Value = -3310.0002
$ws.Range("N70").Value = -3339.9999
$ws.Range("H73").Value = 1063.3334
$ws.Range("I73").Value = 1193.3334
$ws.Range("J73").Value = 933.3333
$ws.Range("K73").Value = 3580.0002
$ws.Range("L73").Value = 2799.9999
$ws.Range("M73").Value = -2644.0002
$ws.Range("N73").Value = -4671.9999
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H138").Value = 3948.8809
$ws.Range("J138").Value = 3940.361
$ws.Range("L138").Value = 11821.083
$ws.Range("N138").Value = -22101.083

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1034.2609
$ws.Range("I2").Value = 982.6667
$ws.Range("K2").Value = 982.6667
$ws.Range("M2").Value = -869.6667
$ws.Range("H116").Value = 1034.2609
$ws.Range("I116").Value = 982.6667
$ws.Range("K116").Value = 982.6667
$ws.Range("M116").Value = 1311.3333
$ws.Range("H122").Value = 1906.6
$ws.Range("I122").Value = 1804.4857
$ws.Range("K122").Value = 5413.4571
$ws.Range("M122").Value = -2963.4571

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1034.2609
$ws.Range("I3").Value = 982.6667
$ws.Range("K3").Value = 982.6667
$ws.Range("M3").Value = -868.6667
$ws.Range("H107").Value = 1945.8387
$ws.Range("I107").Value = 1637.125
$ws.Range("J107").Value = 3004.2856
$ws.Range("K107").Value = 1637.125
$ws.Range("L107").Value = 3004.2856
$ws.Range("M107").Value = 282.875
$ws.Range("N107").Value = -6844.2856

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1601.7142
$ws.Range("I16").Value = 1444.4
$ws.Range("K16").Value = 1444.4
$ws.Range("M16").Value = -1157.4
$ws.Range("H41").Value = 38000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 38000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 38000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -38856
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H105").Value = 3730
$ws.Range("I105").Value = 3973.3333
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 3973.3333
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -2226.3333
$ws.Range("N105").Value = -6494
$ws.Range("H113").Value = 1601.7142
$ws.Range("I113").Value = 1444.4
$ws.Range("K113").Value = 1444.4
$ws.Range("M113").Value = 725.5999999999999
$ws.Range("H134").Value = 76924424
$ws.Range("I134").Value = 83334376
$ws.Range("K134").Value = 250003128
$ws.Range("M134").Value = -250000593

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H112").Value = 125001430
$ws.Range("I112").Value = 1095
$ws.Range("J112").Value = 500002460
$ws.Range("K112").Value = 3285
$ws.Range("L112").Value = 1500007380
$ws.Range("M112").Value = -2177
$ws.Range("N112").Value = -1500009596
$ws.Range("H121").Value = 968.1786
$ws.Range("I121").Value = 264.75
$ws.Range("J121").Value = 1085.4166
$ws.Range("K121").Value = 794.25
$ws.Range("L121").Value = 3256.2498
$ws.Range("M121").Value = 515.75
$ws.Range("N121").Value = -5876.2498
$ws.Range("H122").Value = 1249.2727
$ws.Range("J122").Value = 1266.4375
$ws.Range("L122").Value = 11397.9375
$ws.Range("N122").Value = -16297.9375
$ws.Range("H125").Value = 5000
$ws.Range("J125").Value = 5000
$ws.Range("L125").Value = 15000
$ws.Range("N125").Value = -24840
$ws.Range("H129").Value = 237326.6
$ws.Range("J129").Value = 288063.72
$ws.Range("L129").Value = 864191.1599999999
$ws.Range("N129").Value = -874191.1599999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6643.2354
$ws.Range("I113").Value = 11893.125
$ws.Range("J113").Value = 1976.6666
$ws.Range("K113").Value = 11893.125
$ws.Range("L113").Value = 1976.6666
$ws.Range("M113").Value = -9723.125
$ws.Range("N113").Value = -6316.6666
$ws.Range("H132").Value = 4254437.5
$ws.Range("I132").Value = 7060506
$ws.Range("J132").Value = 45335.5
$ws.Range("K132").Value = 21181518
$ws.Range("L132").Value = 136006.5
$ws.Range("M132").Value = -21178988
$ws.Range("N132").Value = -141066.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3935.389
$ws.Range("I61").Value = 2224.7693
$ws.Range("K61").Value = 2224.7693
$ws.Range("M61").Value = -2022.7693
$ws.Range("H113").Value = 3935.389
$ws.Range("I113").Value = 2224.7693
$ws.Range("K113").Value = 2224.7693
$ws.Range("M113").Value = -54.76929999999993

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4522.222
$ws.Range("I62").Value = 3966.6667
$ws.Range("K62").Value = 3966.6667
$ws.Range("M62").Value = -3342.6667
$ws.Range("H65").Value = 4522.222
$ws.Range("I65").Value = 3966.6667
$ws.Range("K65").Value = 19833.3335
$ws.Range("M65").Value = -16713.3335
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360
$ws.Range("H122").Value = 2110
$ws.Range("I122").Value = 2062.5
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 6187.5
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -3737.5
$ws.Range("N122").Value = -11800
$ws.Range("H136").Value = 33302056
$ws.Range("I136").Value = 44882504
$ws.Range("J136").Value = 8263.125
$ws.Range("K136").Value = 134647512
$ws.Range("L136").Value = 24789.375
$ws.Range("M136").Value = -134644962
$ws.Range("N136").Value = -29889.375
